# Apply the Alvearie -> LinuxForHealth rebrand changes to the FHIR
# quality-measure-population-type ValueSet workbook.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet ---------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/quality-measure-population-type"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date (stored as plain text, not a real date, so force text assignment)
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# Description
$meta.Range("B11").Value = "LinuxForHealth measure population types for qualiity measures"

# --- "Include from Measure Populati" sheet ------------------------------
$include2 = $wb.Worksheets.Item("Include from Measure Populati")

# CodeSystem URL used as the "system" value for this particular include block
$include2.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/measure-population-type"
